$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.04433703455491324
$ws.Range("J2").Value = 0.04433703455491323
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.5699070000000001
$ws.Range("N2").Value = 1.709721
$ws.Range("O2").Value = 0.1475417602688563
$ws.Range("P2").Value = 0.1475417602688563
$ws.Range("Q2").Value = 0.052408077813
$ws.Range("R2").Value = 0.471672700317
$ws.Range("S2").Value = 0.006541564123333009
$ws.Range("T2").Value = 0.006541564123333007
$ws.Range("I3").Value = 0.04433703455491324
$ws.Range("J3").Value = 0.04433703455491323
$ws.Range("N3").Value = 5.084895
$ws.Range("O3").Value = 0.4388051378454766
$ws.Range("P3").Value = 0.4388051378454766
$ws.Range("S3").Value = 0.01945531855952836
$ws.Range("T3").Value = 0.01945531855952836
$ws.Range("I4").Value = 0.04433703455491324
$ws.Range("J4").Value = 0.04433703455491323
$ws.Range("M4").Value = 0.125128
$ws.Range("N4").Value = 0.375384
$ws.Range("O4").Value = 0.03239406671425592
$ws.Range("P4").Value = 0.03239406671425593
$ws.Range("Q4").Value = 0.011506645752
$ws.Range("R4").Value = 0.103559811768
$ws.Range("S4").Value = 0.00143625685528413
$ws.Range("T4").Value = 0.00143625685528413
$ws.Range("I5").Value = 0.04433703455491324
$ws.Range("J5").Value = 0.04433703455491323
$ws.Range("M5").Value = 1.472682666666667
$ws.Range("N5").Value = 4.418048
$ws.Range("O5").Value = 0.3812590351714111
$ws.Range("P5").Value = 0.3812590351714111
$ws.Range("Q5").Value = 0.135426425344
$ws.Range("R5").Value = 1.218837828096
$ws.Range("S5").Value = 0.01690389501676774
$ws.Range("T5").Value = 0.01690389501676773
$ws.Range("G6").Value = 1.982131
$ws.Range("H6").Value = 5.946393
$ws.Range("I6").Value = 0.9556629654450868
$ws.Range("J6").Value = 0.9556629654450867
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.5699070000000001
$ws.Range("N6").Value = 1.709721
$ws.Range("O6").Value = 0.1475417602688563
$ws.Range("P6").Value = 0.1475417602688563
$ws.Range("Q6").Value = 1.129630331817
$ws.Range("R6").Value = 10.166672986353
$ws.Range("S6").Value = 0.1410001961455233
$ws.Range("T6").Value = 0.1410001961455233
$ws.Range("G7").Value = 1.982131
$ws.Range("H7").Value = 5.946393
$ws.Range("I7").Value = 0.9556629654450868
$ws.Range("J7").Value = 0.9556629654450867
$ws.Range("N7").Value = 5.084895
$ws.Range("O7").Value = 0.4388051378454766
$ws.Range("P7").Value = 0.4388051378454766
$ws.Range("Q7").Value = 3.359642670415
$ws.Range("R7").Value = 30.23678403373501
$ws.Range("S7").Value = 0.4193498192859483
$ws.Range("T7").Value = 0.4193498192859482
$ws.Range("G8").Value = 1.982131
$ws.Range("H8").Value = 5.946393
$ws.Range("I8").Value = 0.9556629654450868
$ws.Range("J8").Value = 0.9556629654450867
$ws.Range("M8").Value = 0.125128
$ws.Range("N8").Value = 0.375384
$ws.Range("O8").Value = 0.03239406671425592
$ws.Range("P8").Value = 0.03239406671425593
$ws.Range("Q8").Value = 0.248020087768
$ws.Range("R8").Value = 2.232180789912
$ws.Range("S8").Value = 0.0309578098589718
$ws.Range("T8").Value = 0.0309578098589718
$ws.Range("G9").Value = 1.982131
$ws.Range("H9").Value = 5.946393
$ws.Range("I9").Value = 0.9556629654450868
$ws.Range("J9").Value = 0.9556629654450867
$ws.Range("M9").Value = 1.472682666666667
$ws.Range("N9").Value = 4.418048
$ws.Range("O9").Value = 0.3812590351714111
$ws.Range("P9").Value = 0.3812590351714111
$ws.Range("Q9").Value = 2.919049966762667
$ws.Range("R9").Value = 26.271449700864
$ws.Range("S9").Value = 0.3643551401546433
$ws.Range("T9").Value = 0.3643551401546433